$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "('Avatar', ['Token Creature — Avatar', 'This creature’s power and toughness are each equal to your life total.', '*/*'])"
$ws.Range("A3").Value = "('Beast', ['Token Creature — Beast', '3/3'])"
$ws.Range("A4").Value = "('Elemental', ['Token Creature — Elemental', 'Flying', '4/4'])"
$ws.Range("A5").Value = "('Elemental Shaman', ['Token Creature — Elemental Shaman', '3/1'])"
$ws.Range("A6").Value = "('Elf Warrior', ['Token Creature — Elf Warrior', '1/1'])"
$ws.Range("A7").Value = "('Goblin Rogue', ['Token Creature — Goblin Rogue', '1/1'])"
$ws.Range("A8").Value = "('Kithkin Soldier', ['Token Creature — Kithkin Soldier', '1/1'])"
$ws.Range("A9").Value = "('Merfolk Wizard', ['Token Creature — Merfolk Wizard', '1/1'])"
$ws.Range("A10").Value = "('Shapeshifter', ['Token Creature — Shapeshifter', 'Changeling', '1/1'])"
$ws.Range("A11").Value = "('Wolf', ['Token Creature — Wolf', '2/2'])"

# Delete rows 12 through 37 (the old remaining rows that are no longer needed)
$ws.Range("A12:A37").EntireRow.Delete()
